$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 308. This shifts the existing rows 308-321
# down to 309-322, matching the rest of the diff automatically.
$ws.Rows.Item(308).Insert()

# Populate the newly inserted row 308 with the new data point.
$ws.Range("A308").Value = 10
$ws.Range("B308").Value = "Vega Modelo de Temuco"
$ws.Range("C308").Value = "La Araucanía"
$ws.Range("D308").Value = 44753
$ws.Range("E308").Value = 9
$ws.Range("F308").Value = 100114013
$ws.Range("G308").Value = "Zanahoria"
$ws.Range("H308").Value = "Sin especificar"
$ws.Range("I308").Value = "Primera"
$ws.Range("J308").Value = 100
$ws.Range("K308").Value = 8000
$ws.Range("L308").Value = 8000
$ws.Range("M308").Value = 8000
$ws.Range("N308").Value = "$/saco 25 kilos"
$ws.Range("O308").Value = "Región de La Araucanía"
$ws.Range("P308").Value = 320
$ws.Range("Q308").Value = 25
$ws.Range("R308").Value = "Hortaliza"
